$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data so stale shared strings are garbage-collected on save
$ws.Range("A1:B60").ClearContents()

$ws.Cells.Item(1,1).Value = 'Cluster Name'
$ws.Cells.Item(1,2).Value = 'Activecases'
$ws.Cells.Item(2,1).Value = '3398 BlueCross Elly Kay Mordiallo'
$ws.Cells.Item(2,2).Value = 31
$ws.Cells.Item(3,1).Value = '3564 Waverley Valley Aged Care Glen Waverley'
$ws.Cells.Item(3,2).Value = 16
$ws.Cells.Item(4,1).Value = '3601 Baptcare Westhaven community'
$ws.Cells.Item(4,2).Value = 13
$ws.Cells.Item(5,1).Value = '3646 Mornington Bay Care Community Mount Martha'
$ws.Cells.Item(5,2).Value = 16
$ws.Cells.Item(6,1).Value = '3647 Aurrum Aged Care Reservoir'
$ws.Cells.Item(6,2).Value = 11
$ws.Cells.Item(7,1).Value = '3653 Fronditha Thalpori St Albans Aged Care'
$ws.Cells.Item(7,2).Value = 20
$ws.Cells.Item(8,1).Value = '3975 Aurrum Aged Care Brunswick West'
$ws.Cells.Item(8,2).Value = 13
$ws.Cells.Item(9,1).Value = '4257 BlueCross The Gables Camberwell'
$ws.Cells.Item(9,2).Value = 16
$ws.Cells.Item(10,1).Value = '4295 Hope Aged Care Sunshine West'
$ws.Cells.Item(10,2).Value = 15
$ws.Cells.Item(11,1).Value = '4314 Estia Health Ardeer'
$ws.Cells.Item(11,2).Value = 17
$ws.Cells.Item(12,1).Value = '44095 Myrniong Primary School'
$ws.Cells.Item(12,2).Value = 13
$ws.Cells.Item(13,1).Value = '44304 Brighton Primary School Brighton'
$ws.Cells.Item(13,2).Value = 29
$ws.Cells.Item(14,1).Value = '44414 Buangor Primary School Buangor'
$ws.Cells.Item(14,2).Value = 24
$ws.Cells.Item(15,1).Value = '44490 Armadale Primary School Armadale'
$ws.Cells.Item(15,2).Value = 28
$ws.Cells.Item(16,1).Value = '44593 Torquay P-6 College Torquay'
$ws.Cells.Item(16,2).Value = 29
$ws.Cells.Item(17,1).Value = '44620 Canterbury Primary School Canterbury'
$ws.Cells.Item(17,2).Value = 16
$ws.Cells.Item(18,1).Value = '44623 Brunswick North Primary School Brunswick West'
$ws.Cells.Item(18,2).Value = 27
$ws.Cells.Item(19,1).Value = '44745 Briar Hill Primary School Briar Hill'
$ws.Cells.Item(19,2).Value = 22
$ws.Cells.Item(20,1).Value = '44761 Coburg North Primary School Coburg'
$ws.Cells.Item(20,2).Value = 13
$ws.Cells.Item(21,1).Value = '44799 Eastwood Primary School Ringwood East'
$ws.Cells.Item(21,2).Value = 27
$ws.Cells.Item(22,1).Value = '44828 Cheltenham East Primary School Cheltenham'
$ws.Cells.Item(22,2).Value = 11
$ws.Cells.Item(23,1).Value = '45013 Gladstone Views Primary School'
$ws.Cells.Item(23,2).Value = 20
$ws.Cells.Item(24,1).Value = '45087 Milgate Primary School Doncaster East'
$ws.Cells.Item(24,2).Value = 12
$ws.Cells.Item(25,1).Value = '45147 Maramba Primary School Narre Warren'
$ws.Cells.Item(25,2).Value = 13
$ws.Cells.Item(26,1).Value = '45226 Glen Waverley South Primary School Glen Waverley'
$ws.Cells.Item(26,2).Value = 11
$ws.Cells.Item(27,1).Value = '45257 Roxburgh Rise Primary School Roxburgh Park'
$ws.Cells.Item(27,2).Value = 15
$ws.Cells.Item(28,1).Value = '45305 Lockington Consolidated School'
$ws.Cells.Item(28,2).Value = 10
$ws.Cells.Item(29,1).Value = '45719 St Joseph''s Primary School Numurkah'
$ws.Cells.Item(29,2).Value = 14
$ws.Cells.Item(30,1).Value = '4574 Village Glen Aged Care Residences Mornington'
$ws.Cells.Item(30,2).Value = 10
$ws.Cells.Item(31,1).Value = '45764 Our Lady Help of Christian''s Primary School Brunswick East'
$ws.Cells.Item(31,2).Value = 16
$ws.Cells.Item(32,1).Value = '45858 St Bernard''s Primary Coburg'
$ws.Cells.Item(32,2).Value = 24
$ws.Cells.Item(33,1).Value = '45861 St Oliver Plunkett Primary School Pascoe Vale'
$ws.Cells.Item(33,2).Value = 11
$ws.Cells.Item(34,1).Value = '45975 St Thomas More Primary School Hadfield'
$ws.Cells.Item(34,2).Value = 12
$ws.Cells.Item(35,1).Value = '46074 St Justin''s Catholic Primary School Wheelers Hill'
$ws.Cells.Item(35,2).Value = 12
$ws.Cells.Item(36,1).Value = '46078 Corpus Christi Primary School Werribee'
$ws.Cells.Item(36,2).Value = 34
$ws.Cells.Item(37,1).Value = '46098 Good Samaritan Catholic Primary School Roxburgh Park'
$ws.Cells.Item(37,2).Value = 10
$ws.Cells.Item(38,1).Value = '46101 Emmaus Catholic Primary School Sydenham'
$ws.Cells.Item(38,2).Value = 10
$ws.Cells.Item(39,1).Value = '46135 Wesley College Junior School St Kilda Road'
$ws.Cells.Item(39,2).Value = 10
$ws.Cells.Item(40,1).Value = '46208 Mount Scopus Memorial College Gandel Campus Burwood'
$ws.Cells.Item(40,2).Value = 13
$ws.Cells.Item(41,1).Value = '46306 King''s College Warrnambool'
$ws.Cells.Item(41,2).Value = 12
$ws.Cells.Item(42,1).Value = '46327 Victory Christian College Strathdale'
$ws.Cells.Item(42,2).Value = 10
$ws.Cells.Item(43,1).Value = '50279 Dallas Brooks Community Primary School'
$ws.Cells.Item(43,2).Value = 12
$ws.Cells.Item(44,1).Value = '51529 Sirius College Primary School Dallas'
$ws.Cells.Item(44,2).Value = 11
$ws.Cells.Item(45,1).Value = 'Berwick Chase Primary School Berwick'
$ws.Cells.Item(45,2).Value = 11
$ws.Cells.Item(46,1).Value = 'Brandon Park Primary School Wheelers Hill'
$ws.Cells.Item(46,2).Value = 14
$ws.Cells.Item(47,1).Value = 'Confirmed Omicron Sircuit Bar Fitzroy'
$ws.Cells.Item(47,2).Value = 12
$ws.Cells.Item(48,1).Value = 'Confirmed Omicron Variant The Peel Hotel'
$ws.Cells.Item(48,2).Value = 17
$ws.Cells.Item(49,1).Value = 'Gladstone Views Primary School Gladstone Park'
$ws.Cells.Item(49,2).Value = 13
$ws.Cells.Item(50,1).Value = 'JBS Australia Brooklyn'
$ws.Cells.Item(50,2).Value = 30
$ws.Cells.Item(51,1).Value = 'Kororoit Creek Primary School Burnside Heights'
$ws.Cells.Item(51,2).Value = 32
$ws.Cells.Item(52,1).Value = 'Mambourin Enterprises Allara Deer Park'
$ws.Cells.Item(52,2).Value = 30
$ws.Cells.Item(53,1).Value = 'Mount View Primary School Glen Waverley'
$ws.Cells.Item(53,2).Value = 10
$ws.Cells.Item(54,1).Value = 'Oakleigh South Primary School Oakleigh South'
$ws.Cells.Item(54,2).Value = 11
$ws.Cells.Item(55,1).Value = 'PGL Camp Rumbug Foster North'
$ws.Cells.Item(55,2).Value = 47
$ws.Cells.Item(56,1).Value = 'Rosebud Primary School Rosebud'
$ws.Cells.Item(56,2).Value = 25
$ws.Cells.Item(57,1).Value = 'St Christophers Primary School Airport West'
$ws.Cells.Item(57,2).Value = 12
$ws.Cells.Item(58,1).Value = 'St Louis de Montfort''s School AspendaleOutbreak'
$ws.Cells.Item(58,2).Value = 12
$ws.Cells.Item(59,1).Value = 'St Mary''s Parish Primary School'
$ws.Cells.Item(59,2).Value = 11
$ws.Cells.Item(60,1).Value = 'StarTrack Tullamarine'
$ws.Cells.Item(60,2).Value = 21
$ws.Cells.Item(61,1).Value = 'Thomastown West Primary School Camp Doxa''s Malmsbury'
$ws.Cells.Item(61,2).Value = 11
$ws.Cells.Item(62,1).Value = 'Torquay Hotel Torquay'
$ws.Cells.Item(62,2).Value = 17
